# "code added till purchase journey"
#
# Changes applied to the "signup" sheet (sheet1.xml):
#   1. B2 (mobile)       "11500000222"             -> "11800000222"
#   2. D2 (parent email) "arunpandey14@yopmail.com" -> "arunpandey17@yopmail.com"
#      (the mailto hyperlink on D2 must follow the new address)
#   3. H2 (state)        "- Any -"                  -> "Delhi"
#   4. A new column is inserted before the old column I ("center"),
#      pushing the old I:Q range one column to the right (J:R).
#      New column gets header "center" (I1) and value
#      "Delhi - Pitampura (Engineering)" (I2).
#   5. The hyperlink that used to live on M2 (parent email, 2nd table) must
#      now live on N2 because of the column insertion.
#   6. The sheet view no longer scrolls (topLeftCell is gone) and the
#      active selection becomes E10 instead of E6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1-3. simple value edits -------------------------------------------
$ws.Range("H2").Value = "Delhi"
$ws.Range("B2").Value = "11800000222"
$ws.Range("D2").Value = "arunpandey17@yopmail.com"

# --- 4. insert the new "center" column at I -----------------------------
$ws.Range("I1").EntireColumn.Insert()
$ws.Range("I1").Value = "center"
$ws.Range("I2").Value = "Delhi - Pitampura (Engineering)"

# --- 2/5. rebuild the hyperlinks ----------------------------------------
# Inserting a column does not automatically re-anchor existing
# hyperlinks in this runtime, so remove the old ones and recreate them
# pointing at the (possibly shifted) target cells.
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:arunpandey17@yopmail.com")
$ws.Hyperlinks.Add($ws.Range("N2"), "mailto:arunpandey.pus@aesl.in")

# Adding a hyperlink re-styles the cell with a generic "hyperlink" look;
# restore the original cell formatting (themed underline font on a
# shaded/bordered cell) by copying it back from a cell that still has it.
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("N2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 6. sheet view: drop the manual scroll, move the selection ----------
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("E10").Select()
